# "Add files via upload" -- adds a new worksheet "10_" (1_3_world_pop.xlsx,
# question about world-population carrying capacity) and tweaks the saved
# selection on sheet "1_".

$wb = $excel.ActiveWorkbook

# --- 1) sheet "1_" (xl/worksheets/sheet2.xml): saved selection moves from
#        the stray cell C5 to the actual data block A2:C3. Do this BEFORE
#        adding/activating the new sheet so "1_" does not end up as the
#        active tab. ---
$ws1 = $wb.Worksheets.Item("1_")
$ws1.Range("A2:C3").Select() | Out-Null

# --- 2) Append a brand-new worksheet named "10_" after "9_" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "10_"

# Question / leeway / comment table, matching the other question sheets
$newSheet.Range("A1").Value = "What appears to be carry capacity for world population, based on the US and UN data?  Adjust our parameters to better match their data to find an answer."
$newSheet.Range("B1").Value = "Leeway"
$newSheet.Range("C1").Value = "Comment"
$newSheet.Range("A2").Value = 11.8
$newSheet.Range("B2").Value = 0.3
$newSheet.Range("C3").Value = "Adjust 'r' and 'K' to match the official data.  Which represents the carrying capacity?"

# Wrap text over the whole used block (also stamps the wrap-text style onto
# the otherwise-empty cells, same as the rest of the workbook's sheets)
$newSheet.Range("A1:C5").WrapText = $true

# Row heights used by the author for the wrapped header row / comment row
$newSheet.Rows.Item(1).RowHeight = 90
$newSheet.Rows.Item(3).RowHeight = 60

# Column widths matching the other question sheets' proportions
$newSheet.Columns.Item(1).ColumnWidth = 26.6
$newSheet.Columns.Item(2).ColumnWidth = 32.92
$newSheet.Columns.Item(3).ColumnWidth = 26.6

# Leave the new sheet as the active tab with its last-used selection --
# must be the final UI action so "10_" becomes tabSelected / activeTab.
$newSheet.Range("B11").Select()

Write-Output "Added sheet '10_' and updated selection on '1_'"
